$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 18561
$ws.Range("E3").Value = 15392
$ws.Range("E4").Value = 7834
$ws.Range("E5").Value = 13540
$ws.Range("E6").Value = 18132
$ws.Range("E7").Value = 4648
$ws.Range("E8").Value = 13210
$ws.Range("E9").Value = 11792
$ws.Range("E10").Value = 6013
$ws.Range("E11").Value = 16243
$ws.Range("E12").Value = 2275
$ws.Range("E13").Value = 1979
